$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 held only the "Docentes responsáveis:" value (B13/C13 = "519033 -
# Carlos Yujiro Shigue") with no A13 label. That row is removed entirely, which
# shifts every subsequent row up by one (row heights / A-column labels then line
# up exactly with what used to be one row below).
$ws.Rows(13).Delete()

# After the shift, several B/C cells need their text content replaced.
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2016"
$ws.Range("C15").Value = "01/01/2016"

$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("B19").Value = "As atividades práticas e os projetos que serão desenvolvidos durante as aulas serão avaliados por docentes e pelos alunos (processo de avaliação crítica)."
$ws.Range("C19").Value = "As atividades práticas e os projetos que serão desenvolvidos durante as aulas serão avaliados por docentes e pelos alunos (processo de avaliação crítica)."

$ws.Range("B20").Value = "A média final será uma composição de fatores relativos à participação do aluno nos trabalhos desenvolvidos, conjuntamente com o rendimento de seu grupo."
$ws.Range("C20").Value = "A média final será uma composição de fatores relativos à participação do aluno nos trabalhos desenvolvidos, conjuntamente com o rendimento de seu grupo."

$ws.Range("B21").Value = "Devido às características da disciplina, não será oferecida recuperação."
$ws.Range("C21").Value = "Devido às características da disciplina, não será oferecida recuperação."
